$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# shape id=52 (index 1)
$sh1 = $s.Shapes.Item(1)
$sh1.Left = 492.7032470703125
$sh1.Top = 204.84796142578125

# shape id=68 (index 2)
$sh2 = $s.Shapes.Item(2)
$sh2.Left = 185.5648956298828
$sh2.Top = 94.68315124511719

# shape id=58 (index 3)
$sh3 = $s.Shapes.Item(3)
$sh3.Left = 492.93560791015625
$sh3.Top = 94.49385833740234

# shape id=59 (index 4)
$sh4 = $s.Shapes.Item(4)
$sh4.Left = 185.9454345703125
$sh4.Top = 205.47071838378906

# shape id=38 (index 7)
$sh7 = $s.Shapes.Item(7)
$sh7.Left = 177.48654174804688
$sh7.Top = 71.21732330322266
$sh7.Width = 193.0779571533203
$sh7.Height = 444.70843505859375

# shape id=42 (index 8)
$sh8 = $s.Shapes.Item(8)
$sh8.Left = 109.7142562866211
$sh8.Top = 62.946929931640625
$sh8.Width = 576.9039916992188
$sh8.Height = 460.6181945800781

# shape id=45 (index 9)
$sh9 = $s.Shapes.Item(9)
$sh9.Left = 109.92512512207031
$sh9.Top = 63.07307434082031

# shape id=60 (index 10)
$sh10 = $s.Shapes.Item(10)
$sh10.Left = 185.9454345703125
$sh10.Top = 204.93023681640625

# shape id=69 (index 11)
$sh11 = $s.Shapes.Item(11)
$sh11.Left = 185.30284118652344
$sh11.Top = 94.68890380859375

# shape id=67 (index 12)
$sh12 = $s.Shapes.Item(12)
$sh12.Left = 484.09222412109375
$sh12.Top = 71.21732330322266
$sh12.Width = 191.7878875732422
$sh12.Height = 444.3194580078125

# shape id=79 (index 13)
$sh13 = $s.Shapes.Item(13)
$sh13.Left = 492.7032470703125
$sh13.Top = 204.93023681640625

# shape id=85 (index 14)
$sh14 = $s.Shapes.Item(14)
$sh14.Left = 493.2381286621094
$sh14.Top = 94.68315124511719

# shape id=93 (index 15)
$sh15 = $s.Shapes.Item(15)
$sh15.Left = 416.08929443359375
$sh15.Top = 118.18118286132812

# shape id=94 (index 16)
$sh16 = $s.Shapes.Item(16)
$sh16.Left = 303.99530029296875
$sh16.Top = 128.9480438232422

# shape id=96 (index 17)
$sh17 = $s.Shapes.Item(17)
$sh17.Left = 280.45947265625
$sh17.Top = 156.1859130859375

# shape id=97 (index 18)
$sh18 = $s.Shapes.Item(18)
$sh18.Left = 491.96575927734375
$sh18.Top = 156.37322998046875

# shape id=98 (index 19)
$sh19 = $s.Shapes.Item(19)
$sh19.Left = 212.24354553222656
$sh19.Top = 125.74866485595703

# shape id=99 (index 20)
$sh20 = $s.Shapes.Item(20)
$sh20.Left = 605.828125
$sh20.Top = 125.34803771972656

# shape id=100 (index 21)
$sh21 = $s.Shapes.Item(21)
$sh21.Left = 170.88134765625
$sh21.Top = 156.0267791748047

# shape id=101 (index 22)
$sh22 = $s.Shapes.Item(22)
$sh22.Left = 574.2802124023438
$sh22.Top = 156.8473358154297

# shape id=73 (index 23)
$sh23 = $s.Shapes.Item(23)
$sh23.Left = 517.6026000976562
$sh23.Top = 128.9480438232422

# shape id=119 (index 31)
$sh31 = $s.Shapes.Item(31)
$sh31.Left = 236.25331115722656
$sh31.Top = 254.44622802734375

# shape id=120 (index 32)
$sh32 = $s.Shapes.Item(32)
$sh32.Left = 374.7177429199219
$sh32.Top = 251.2428436279297

# shape id=121 (index 33)
$sh33 = $s.Shapes.Item(33)
$sh33.Left = 408.56646728515625
$sh33.Top = 217.54685974121094

# shape id=122 (index 34)
$sh34 = $s.Shapes.Item(34)
$sh34.Left = 296.6952209472656
$sh34.Top = 243.40631103515625

# shape id=123 (index 35)
$sh35 = $s.Shapes.Item(35)
$sh35.Left = 259.7890625
$sh35.Top = 278.63922119140625

# shape id=124 (index 36)
$sh36 = $s.Shapes.Item(36)
$sh36.Left = 235.60362243652344
$sh36.Top = 310.4985046386719

# shape id=125 (index 37)
$sh37 = $s.Shapes.Item(37)
$sh37.Left = 557.3072509765625
$sh37.Top = 277.7640380859375

# shape id=126 (index 38)
$sh38 = $s.Shapes.Item(38)
$sh38.Left = 535.21142578125
$sh38.Top = 308.8627014160156

# shape id=127 (index 39)
$sh39 = $s.Shapes.Item(39)
$sh39.Left = 374.6523742675781
$sh39.Top = 319.45245361328125

# shape id=128 (index 40)
$sh40 = $s.Shapes.Item(40)
$sh40.Left = 408.25372314453125
$sh40.Top = 286.963623046875

# shape id=129 (index 41)
$sh41 = $s.Shapes.Item(41)
$sh41.Left = 297.8950500488281
$sh41.Top = 301.5555419921875

# shape id=130 (index 42)
$sh42 = $s.Shapes.Item(42)
$sh42.Left = 449.7269592285156
$sh42.Top = 301.63018798828125

# shape id=131 (index 43)
$sh43 = $s.Shapes.Item(43)
$sh43.Left = 259.7890625
$sh43.Top = 338.2967834472656

# shape id=132 (index 44)
$sh44 = $s.Shapes.Item(44)
$sh44.Left = 236.25331115722656
$sh44.Top = 365.53466796875

# shape id=134 (index 46)
$sh46 = $s.Shapes.Item(46)
$sh46.Left = 533.7715454101562
$sh46.Top = 365.7878112792969

# shape id=136 (index 48)
$sh48 = $s.Shapes.Item(48)
$sh48.Left = 214.60922241210938
$sh48.Top = 419.4653625488281

# shape id=138 (index 50)
$sh50 = $s.Shapes.Item(50)
$sh50.Left = 533.7715454101562
$sh50.Top = 419.718505859375

# shape id=140 (index 52)
$sh52 = $s.Shapes.Item(52)
$sh52.Left = 253.77040100097656
$sh52.Top = 418.5812072753906

# shape id=141 (index 53)
$sh53 = $s.Shapes.Item(53)
$sh53.Left = 259.750244140625
$sh53.Top = 454.92852783203125

# shape id=142 (index 54)
$sh54 = $s.Shapes.Item(54)
$sh54.Left = 233.48370361328125
$sh54.Top = 490.4989013671875

# shape id=143 (index 55)
$sh55 = $s.Shapes.Item(55)
$sh55.Left = 563.3630981445312
$sh55.Top = 455.3800964355469

# shape id=144 (index 56)
$sh56 = $s.Shapes.Item(56)
$sh56.Left = 537.0965576171875
$sh56.Top = 490.95050048828125

# shape id=64 (index 57)
$sh57 = $s.Shapes.Item(57)
$sh57.Left = 283.1827697753906
$sh57.Top = 118.8271713256836
$sh57.Width = 288.8567810058594
$sh57.Height = 68.24984741210938

# shape id=61 (index 58)
$sh58 = $s.Shapes.Item(58)
$sh58.Left = 233.48370361328125
$sh58.Top = 439.9892272949219

# shape id=62 (index 59)
$sh59 = $s.Shapes.Item(59)
$sh59.Left = 415.98260498046875
$sh59.Top = 439.4694519042969

